$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D against Excel auto-converting numeric-looking text
# (e.g. "0.999", "24.40") into actual numbers: force text format first,
# then restore the default "Normal" style once the text values are set.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '57.321.45'
$ws.Range('E2').Value = '  +1.57%  '

$ws.Range('D3').Value = '3.022.94'
$ws.Range('E3').Value = '  +0.91%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = '519.00'
$ws.Range('E5').Value = '  +4.88%  '

$ws.Range('D6').Value = '140.89'
$ws.Range('E6').Value = '  +5.42%  '

$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').Value = '0.438'

$ws.Range('D9').Value = '7.59'
$ws.Range('E9').Value = '  +5.36%  '

$ws.Range('D10').Value = '0.110'
$ws.Range('E10').Value = '  +6.37%  '

$ws.Range('D11').Value = '0.361'
$ws.Range('E11').Value = '  +2.58%  '

$ws.Range('E12').Value = '  +2.33%  '

$ws.Range('D13').Value = '3.535.85'
$ws.Range('E13').Value = '  +0.96%  '

$ws.Range('D14').Value = '26.16'
$ws.Range('E14').Value = '  +5.31%  '

$ws.Range('D15').Value = '0.0000161'
$ws.Range('E15').Value = '  +11.99%  '

$ws.Range('D16').Value = '57.258.27'
$ws.Range('E16').Value = '  +1.57%  '

$ws.Range('D17').Value = '3.015.83'
$ws.Range('E17').Value = '  +0.95%  '

$ws.Range('D18').Value = '6.03'
$ws.Range('E18').Value = '  +3.76%  '

$ws.Range('D19').Value = '12.72'
$ws.Range('E19').Value = '  +3.30%  '

$ws.Range('D20').Value = '7.97'
$ws.Range('E20').Value = '  +3.26%  '

$ws.Range('D21').Value = '331.73'
$ws.Range('E21').Value = '  +3.01%  '

$ws.Range('E22').Value = '  +0.00%  '

$ws.Range('D23').Value = '0.488'
$ws.Range('E23').Value = '  +5.59%  '

$ws.Range('D24').Value = '63.96'
$ws.Range('E24').Value = '  +4.99%  '

$ws.Range('D25').Value = '0.173'
$ws.Range('E25').Value = '  +5.60%  '

$ws.Range('E26').Value = '  +0.42%  '

$ws.Range('D27').Value = '0.0₃0929'
$ws.Range('E27').Value = '  +6.88%  '

$ws.Range('D28').Value = '6.76'
$ws.Range('E28').Value = '  +3.49%  '

$ws.Range('D29').Value = '7.16'
$ws.Range('E29').Value = '  +7.27%  '

$ws.Range('D30').Value = '1.83'
$ws.Range('E30').Value = '  +6.48%  '

$ws.Range('D31').Value = '1.23'
$ws.Range('E31').Value = '  +4.86%  '

$ws.Range('D32').Value = '20.79'
$ws.Range('E32').Value = '  +4.97%  '

$ws.Range('D33').Value = '157.64'
$ws.Range('E33').Value = '  +4.27%  '

$ws.Range('D34').Value = '4.64'
$ws.Range('E34').Value = '  +4.11%  '

$ws.Range('D35').Value = '5.78'
$ws.Range('E35').Value = '  +3.07%  '

$ws.Range('D36').Value = '1.30'
$ws.Range('E36').Value = '  +1.37%  '

$ws.Range('B37').Value = 'EnergySwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D37').Value = '24.40'
$ws.Range('E37').Value = '  +2.81%  '

$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.0681'
$ws.Range('E38').Value = '  +3.19%  '

$ws.Range('D39').Value = '3.051.67'
$ws.Range('E39').Value = '  +0.89%  '

$ws.Range('D40').Value = '37.35'
$ws.Range('E40').Value = '  +1.88%  '

$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.15%  '

$ws.Range('D42').Value = '2.301.84'
$ws.Range('E42').Value = '  +6.17%  '

$ws.Range('D43').Value = '0.652'
$ws.Range('E43').Value = '  +1.93%  '

$ws.Range('D44').Value = '3.74'
$ws.Range('E44').Value = '  +5.58%  '

$ws.Range('D45').Value = '1.44'
$ws.Range('E45').Value = '  +1.95%  '

$ws.Range('E46').Value = '  +0.39%  '

$ws.Range('D47').Value = '2.01'
$ws.Range('E47').Value = '  +9.51%  '

$ws.Range('D48').Value = '0.0242'
$ws.Range('E48').Value = '  +2.26%  '

$ws.Range('D49').Value = '5.91'
$ws.Range('E49').Value = '  +6.24%  '

$ws.Range('D50').Value = '19.52'
$ws.Range('E50').Value = '  +0.94%  '

$ws.Range('D51').Value = '0.0884'
$ws.Range('E51').Value = '  +4.20%  '

# Restore default styling on column D now that the text values are committed
$ws.Range("D2:D51").Style = "Normal"

